$p = $ppt.ActivePresentation

# Slide 2: "Some notes on the second slide."
$s2 = $p.Slides.Item(2)
$notes2 = $s2.NotesPage
$notes2.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Some notes on the second slide."

# Slide 3: "Final notes on the third slide." + "Second line of notes."
$s3 = $p.Slides.Item(3)
$notes3 = $s3.NotesPage
$notes3.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Final notes on the third slide.`rSecond line of notes."
